# Update CHA Yearly Financials sheet with latest reported figures
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CHA")

# Income Statement
$ws.Range("D8").Value = 54352000
$ws.Range("E8").Value = 52319600
$ws.Range("F8").Value = 49153700
$ws.Range("G8").Value = 48143300
$ws.Range("H8").Value = 47726300
$ws.Range("I8").Value = 42010900
$ws.Range("J8").Value = 36370500
$ws.Range("D9").Value = 13734500
$ws.Range("E9").Value = 14290800
$ws.Range("F9").Value = 12245900
$ws.Range("G9").Value = 11399500
$ws.Range("H9").Value = 15995900
$ws.Range("I9").Value = 13678200
$ws.Range("J9").Value = 10189100
$ws.Range("D10").Value = 40617600
$ws.Range("E10").Value = 38028700
$ws.Range("F10").Value = 36907800
$ws.Range("G10").Value = 36743800
$ws.Range("H10").Value = 31730400
$ws.Range("I10").Value = 28332700
$ws.Range("J10").Value = 26181500
$ws.Range("D15").Value = 11123500
$ws.Range("E15").Value = 10083300
$ws.Range("F15").Value = 10042000
$ws.Range("G15").Value = 9846300
$ws.Range("H15").Value = 10252600
$ws.Range("I15").Value = 7369300
$ws.Range("J15").Value = 7603500
$ws.Range("D17").Value = 50312300
$ws.Range("E17").Value = 48279900
$ws.Range("F17").Value = 45229400
$ws.Range("G17").Value = 43912400
$ws.Range("H17").Value = 43649800
$ws.Range("I17").Value = 38866600
$ws.Range("J17").Value = 32789900
$ws.Range("D18").Value = 4039700
$ws.Range("E18").Value = 4039700
$ws.Range("F18").Value = 3924300
$ws.Range("G18").Value = 4230900
$ws.Range("H18").Value = 4076500
$ws.Range("I18").Value = 3144200
$ws.Range("J18").Value = 3580700
$ws.Range("D20").Value = 195800
$ws.Range("E20").Value = 88700
$ws.Range("F20").Value = 715900
$ws.Range("G20").Value = 59200
$ws.Range("H20").Value = 167900
$ws.Range("I20").Value = 112900
$ws.Range("J20").Value = 88300
$ws.Range("D21").Value = 15378400
$ws.Range("E21").Value = 14229400
$ws.Range("F21").Value = 14699800
$ws.Range("G21").Value = 14153600
$ws.Range("H21").Value = 14514900
$ws.Range("I21").Value = 10639400
$ws.Range("J21").Value = "NA"
$ws.Range("D22").Value = 532200
$ws.Range("E22").Value = 549400
$ws.Range("F22").Value = 678700
$ws.Range("G22").Value = 838500
$ws.Range("H22").Value = 817900
$ws.Range("I22").Value = 319700
$ws.Range("J22").Value = 402200
$ws.Range("D23").Value = 3703300
$ws.Range("E23").Value = 3579100
$ws.Range("F23").Value = 3961500
$ws.Range("G23").Value = 3451600
$ws.Range("H23").Value = 3426500
$ws.Range("I23").Value = 2937500
$ws.Range("J23").Value = 3266800
$ws.Range("D24").Value = 919000
$ws.Range("E24").Value = 889400
$ws.Range("F24").Value = 972200
$ws.Range("G24").Value = 816000
$ws.Range("H24").Value = 804700
$ws.Range("I24").Value = 705400
$ws.Range("J24").Value = 803800
$ws.Range("D26").Value = 2784300
$ws.Range("E26").Value = 2689600
$ws.Range("F26").Value = 2989300
$ws.Range("G26").Value = 2635600
$ws.Range("H26").Value = 2621800
$ws.Range("I26").Value = 2232100
$ws.Range("J26").Value = 2463000
$ws.Range("D27").Value = 2762900
$ws.Range("E27").Value = 2674100
$ws.Range("F27").Value = 2976200
$ws.Range("G27").Value = 2623900
$ws.Range("H27").Value = 2603900
$ws.Range("I27").Value = 2215000
$ws.Range("J27").Value = 2448800
$ws.Range("D32").Value = -195800
$ws.Range("E32").Value = -88700
$ws.Range("F32").Value = -715900
$ws.Range("G32").Value = -59200
$ws.Range("H32").Value = -167900
$ws.Range("I32").Value = -112900
$ws.Range("J32").Value = -88300
$ws.Range("D33").Value = 2762900
$ws.Range("E33").Value = 2674100
$ws.Range("F33").Value = 2976200
$ws.Range("G33").Value = 2623900
$ws.Range("H33").Value = 2603900
$ws.Range("I33").Value = 2215000
$ws.Range("J33").Value = 2448800
$ws.Range("D35").Value = 2762900
$ws.Range("E35").Value = 2674100
$ws.Range("F35").Value = 2976200
$ws.Range("G35").Value = 2623900
$ws.Range("H35").Value = 2603900
$ws.Range("I35").Value = 2215000
$ws.Range("J35").Value = 2448800
# Balance Sheet
$ws.Range("D41").Value = 2636200
$ws.Range("E41").Value = 3286800
$ws.Range("F41").Value = 4588200
$ws.Range("G41").Value = 2769300
$ws.Range("H41").Value = 2172600
$ws.Range("I41").Value = 3320700
$ws.Range("J41").Value = 3631600
$ws.Range("D42").Value = 704500
$ws.Range("E42").Value = 860900
$ws.Range("F42").Value = 515300
$ws.Range("G42").Value = 468200
$ws.Range("H42").Value = 551800
$ws.Range("I42").Value = 1534100
$ws.Range("J42").Value = 698400
$ws.Range("D43").Value = 5771500
$ws.Range("E43").Value = 5139700
$ws.Range("F43").Value = 4704300
$ws.Range("G43").Value = 4186600
$ws.Range("H43").Value = 3602200
$ws.Range("I43").Value = 3467700
$ws.Range("J43").Value = 3445600
$ws.Range("D44").Value = 611900
$ws.Range("E44").Value = 757800
$ws.Range("F44").Value = 932200
$ws.Range("G44").Value = 627000
$ws.Range("H44").Value = 968100
$ws.Range("I44").Value = 879800
$ws.Range("J44").Value = 718700
$ws.Range("D45").Value = 894600
$ws.Range("E45").Value = 956900
$ws.Range("F45").Value = 852000
$ws.Range("G45").Value = 785500
$ws.Range("H45").Value = 538900
$ws.Range("I45").Value = 475500
$ws.Range("J45").Value = 347900
$ws.Range("D46").Value = 10618700
$ws.Range("E46").Value = 11002200
$ws.Range("F46").Value = 11592000
$ws.Range("G46").Value = 8836800
$ws.Range("H46").Value = 7833500
$ws.Range("I46").Value = 9677800
$ws.Range("J46").Value = 8842300
$ws.Range("D47").Value = 5473400
$ws.Range("E47").Value = 5358600
$ws.Range("F47").Value = 5357200
$ws.Range("G47").Value = 753600
$ws.Range("H47").Value = 316400
$ws.Range("I47").Value = 242200
$ws.Range("J47").Value = 242400
$ws.Range("D48").Value = 71142300
$ws.Range("E48").Value = 69761200
$ws.Range("F48").Value = 65758100
$ws.Range("G48").Value = 63231100
$ws.Range("H48").Value = 62109300
$ws.Range("I48").Value = 60288100
$ws.Range("J48").Value = 42645900
$ws.Range("D49").Value = 6279400
$ws.Range("E49").Value = 6109600
$ws.Range("F49").Value = 6034200
$ws.Range("G49").Value = 5773300
$ws.Range("H49").Value = 5633900
$ws.Range("I49").Value = 5807600
$ws.Range("J49").Value = 5585100
$ws.Range("D52").Value = 4614100
$ws.Range("E52").Value = 4614500
$ws.Range("F52").Value = 4691700
$ws.Range("G52").Value = 4703900
$ws.Range("H52").Value = 4728900
$ws.Range("I52").Value = 4878400
$ws.Range("J52").Value = 4890400
$ws.Range("D54").Value = 98127800
$ws.Range("E54").Value = 96846100
$ws.Range("F54").Value = 93433100
$ws.Range("G54").Value = 83298700
$ws.Range("H54").Value = 80622100
$ws.Range("I54").Value = 80894100
$ws.Range("J54").Value = 62206100
$ws.Range("D57").Value = 18193900
$ws.Range("E57").Value = 18574100
$ws.Range("F57").Value = 18202000
$ws.Range("G57").Value = 13293500
$ws.Range("H57").Value = 12300400
$ws.Range("I57").Value = 18411400
$ws.Range("J57").Value = 8760500
$ws.Range("D58").Value = 8274600
$ws.Range("E58").Value = 15302300
$ws.Range("F58").Value = 7681400
$ws.Range("G58").Value = 6538600
$ws.Range("H58").Value = 7088100
$ws.Range("I58").Value = 2483600
$ws.Range("J58").Value = 3109600
$ws.Range("D59").Value = 14404800
$ws.Range("E59").Value = 13486200
$ws.Range("F59").Value = 12099000
$ws.Range("G59").Value = 10788500
$ws.Range("H59").Value = 10308100
$ws.Range("I59").Value = 7816500
$ws.Range("J59").Value = 7016800
$ws.Range("D60").Value = 40873300
$ws.Range("E60").Value = 47362500
$ws.Range("F60").Value = 37982400
$ws.Range("G60").Value = 30620700
$ws.Range("H60").Value = 29696500
$ws.Range("I60").Value = 28711500
$ws.Range("J60").Value = 18887000
$ws.Range("D61").Value = 7216000
$ws.Range("E61").Value = 1398000
$ws.Range("F61").Value = 9633400
$ws.Range("G61").Value = 9337700
$ws.Range("H61").Value = 9293000
$ws.Range("I61").Value = 12328900
$ws.Range("J61").Value = 4623000
$ws.Range("D62").Value = 1553400
$ws.Range("E62").Value = 1136400
$ws.Range("F62").Value = 589200
$ws.Range("G62").Value = 285400
$ws.Range("H62").Value = 276000
$ws.Range("I62").Value = 372200
$ws.Range("J62").Value = 568300
$ws.Range("D66").Value = 49765900
$ws.Range("E66").Value = 50041000
$ws.Range("F66").Value = 48348600
$ws.Range("G66").Value = 40381000
$ws.Range("H66").Value = 39402600
$ws.Range("I66").Value = 41555200
$ws.Range("J66").Value = 24195100
$ws.Range("D72").Value = 34886700
$ws.Range("E72").Value = 33291500
$ws.Range("F72").Value = 31599200
$ws.Range("G72").Value = 29451400
$ws.Range("H72").Value = 27753700
$ws.Range("I72").Value = 25854800
$ws.Range("J72").Value = 24526500
$ws.Range("D76").Value = 48361900
$ws.Range("E76").Value = 46805100
$ws.Range("F76").Value = 45084600
$ws.Range("G76").Value = 42917600
$ws.Range("H76").Value = 41219500
$ws.Range("I76").Value = 39338900
$ws.Range("J76").Value = 38011100
# Cash Flow Statement
$ws.Range("D81").Value = 2762900
$ws.Range("E81").Value = 2674100
$ws.Range("F81").Value = 2976200
$ws.Range("G81").Value = 2623900
$ws.Range("H81").Value = 2603900
$ws.Range("I81").Value = 2215000
$ws.Range("J81").Value = 2448800
$ws.Range("D83").Value = 11123500
$ws.Range("E83").Value = 10083300
$ws.Range("F83").Value = 10042000
$ws.Range("G83").Value = 9846300
$ws.Range("H83").Value = 10252600
$ws.Range("I83").Value = 7369300
$ws.Range("J83").Value = "NA"
$ws.Range("D89").Value = 14321900
$ws.Range("E89").Value = 15009400
$ws.Range("F89").Value = 16139600
$ws.Range("G89").Value = 14307500
$ws.Range("H89").Value = 13112200
$ws.Range("I89").Value = 10487700
$ws.Range("J89").Value = 10835300
$ws.Range("D91").Value = -12961200
$ws.Range("E91").Value = -14348000
$ws.Range("F91").Value = -15122700
$ws.Range("G91").Value = -11913300
$ws.Range("H91").Value = -10525400
$ws.Range("I91").Value = -7424700
$ws.Range("J91").Value = -7197100
$ws.Range("D94").Value = -12653900
$ws.Range("E94").Value = -14699000
$ws.Range("F94").Value = -15174900
$ws.Range("G94").Value = -12126300
$ws.Range("H94").Value = -16020600
$ws.Range("I94").Value = -7161100
$ws.Range("J94").Value = "NA"
$ws.Range("D100").Value = -2396400
$ws.Range("E100").Value = -1418100
$ws.Range("F100").Value = 713700
$ws.Range("G100").Value = -1532600
$ws.Range("H100").Value = 836600
$ws.Range("I100").Value = -2938800
$ws.Range("J100").Value = "NA"
$ws.Range("D101").Value = -44400
$ws.Range("E101").Value = 31300
$ws.Range("F101").Value = 18400
$ws.Range("G101").Value = -600
$ws.Range("H101").Value = -10200
$ws.Range("I101").Value = -400
$ws.Range("J101").Value = "NA"
$ws.Range("D102").Value = -772800
$ws.Range("E102").Value = -1076300
$ws.Range("F102").Value = 1696800
$ws.Range("G102").Value = 648000
$ws.Range("H102").Value = -2082000
$ws.Range("I102").Value = 387400
$ws.Range("J102").Value = 229700
